# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (total) sheet
#    and populate it with the per-fund holding detail for the new quarter.
# 2. Insert a new row at the top of the "总计" sheet's data (row 2) with the
#    2022-Q1 summary figures, pushing the existing quarterly rows down and
#    renumbering the leading index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: new "2022-Q1" sheet, inserted before "总计"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($total)
$newSheet.Name = "2022-Q1"

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'001092"
$newSheet.Range("C2").Value = "广发纳斯达克生物科技指数(QDII)（人民币）"
$newSheet.Range("D2").Value = "'1.34"
$newSheet.Range("E2").Value = "'82.00"
$newSheet.Range("F2").Value = "'5.10"
$newSheet.Range("G2").Value = "'0.0683"
$newSheet.Range("H2").Value = 4

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'001093"
$newSheet.Range("C3").Value = "广发纳斯达克生物科技指数(QDII)（美元）"
$newSheet.Range("D3").Value = "'1.34"
$newSheet.Range("E3").Value = "'82.00"
$newSheet.Range("F3").Value = "'5.10"
$newSheet.Range("G3").Value = "'0.0683"
$newSheet.Range("H3").Value = 4

# ---------------------------------------------------------------------
# Step 2: prepend the 2022-Q1 row to the "总计" summary sheet
# ---------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")
$totalWs.Rows.Item(2).Insert()

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("C2").Value = 2
$totalWs.Range("D2").Value = 0.14

# Renumber the index column (A) for the rows that shifted down one place.
for ($r = 3; $r -le 7; $r++) {
  $totalWs.Cells.Item($r, 1).Value = $r - 2
}

# The inserted row's index cell (A2) doesn't inherit the bold/bordered
# "index column" look the other rows have (row 1 has no cell above A2 to
# copy from), so reapply it to match A3:A7.
$a2 = $totalWs.Range("A2")
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1
